# Swap the contents of rows 35 and 36 on the active sheet.
# (The two records were reordered; row 35 now holds what used to be in
# row 36 and vice versa, including the substrate columns AJ/AK/AO that
# only applied to the "gran" / Picea abies record.)
#
# Only the columns whose values actually differ between the two rows are
# touched, so columns that already hold identical text (e.g. the
# date-like strings in Y/AA) are left completely alone and can't get
# reinterpreted as a different cell type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row1 = 35
$row2 = 36

# Columns A,B,E,F,G,H,Q,R,AJ,AK,AO
$cols = @(1, 2, 5, 6, 7, 8, 17, 18, 36, 37, 41)

foreach ($c in $cols) {
    $cell1 = $ws.Cells.Item($row1, $c)
    $cell2 = $ws.Cells.Item($row2, $c)

    $v1 = $cell1.Value()
    $v2 = $cell2.Value()

    if ($v2 -eq $null -or $v2 -eq "") {
        $cell1.ClearContents()
    } else {
        $cell1.Value() = $v2
    }

    if ($v1 -eq $null -or $v1 -eq "") {
        $cell2.ClearContents()
    } else {
        $cell2.Value() = $v1
    }
}
